$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FindAndBookFlight")

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "selectFlight"

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "bookFlight"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E7"))

$ws.Range("B7").Select()
